$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B93: convert stored text "1" into a true numeric value 1
$ws.Range("B93").Value = 1

# Append new row 94 with the new annotation record
$ws.Range("A94").Value = "Ying Tang"

# B94 must remain TEXT "5" (not numeric) - use a leading apostrophe to force
# text entry, then reset the style pointer so no stray NumberFormat sticks.
$ws.Range("B94").Value = "'5"
$ws.Range("B94").Style = "Normal"

$ws.Range("C94").Value = "We are also happy ,to be exciting"
$ws.Range("D94").Value = "APC"
$ws.Range("E94").Value = "OTH"
$ws.Range("F94").Value = "e9624372-e81d-40ef-b27a-4327fdc73888"
$ws.Range("G94").Value = "BkN_r2lR-_annotated.xlsx"
$ws.Range("H94").Value = "We are also happy that AnonReviewer2 found the list of possible applications, provided in reply to the challenge posted in the review, to be exciting."
